$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 222; this shifts the existing rows
# 222-358 down to 223-359 and extends the used range to row 359.
$ws.Rows("222:222").Insert()

# Populate the newly inserted row 222 with its data.
$ws.Range("A222").Value = 5
$ws.Range("B222").Value = "Macroferia Regional de Talca"
$ws.Range("C222").Value = "Maule"
$ws.Range("D222").Value = 44777
$ws.Range("E222").Value = 7
$ws.Range("F222").Value = 100114014
$ws.Range("G222").Value = "Betarraga"
$ws.Range("H222").Value = "Sin especificar"
$ws.Range("I222").Value = "Primera"
$ws.Range("J222").Value = 3000
$ws.Range("K222").Value = 750
$ws.Range("L222").Value = 750
$ws.Range("M222").Value = 750
$ws.Range("N222").Value = "`$/paquete 5 unidades"
$ws.Range("O222").Value = "Región del Maule"
$ws.Range("P222").Value = 150
$ws.Range("Q222").Value = 5
$ws.Range("R222").Value = "Hortaliza"

# Match the date number format used by the other rows' Fecha column.
$ws.Range("D222").NumberFormat = $ws.Range("D223").NumberFormat
